# Mark mosquito net questions as required.
#
# The "survey" worksheet is an XLSForm-style question sheet with columns:
#   A: clause, B: type, C: values_list, D: name, E: display.prompt
# A new column F ("required") is added, and the three question rows
# (obtain_net_how_long_ago / net_sleep_last_night / net_used_other_purpose)
# are flagged as required (value 1).

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

$survey.Range("F1").Value = "required"
$survey.Range("F2").Value = 1
$survey.Range("F3").Value = 1
$survey.Range("F4").Value = 1

# Bump the form_version on the "settings" sheet, as is customary when a
# form's questions change.
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20210421001

# The workbook was saved with the "survey" sheet active (instead of
# "settings"), with cell F5 selected just below the newly-added column.
$survey.Activate()
$null = $survey.Range("F5").Select()
